$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 48, shifting rows 48:209 down to 49:210.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new September entry.
$ws.Range("R48").Value = "bal axisbank axis"
$ws.Range("S48").Value = "2024-09-25 07:22:34"
